# Paragraph-level OOXML fragments describing the lrc/test.docx playlist content.
# Each entry becomes one Word paragraph; proofErr / lastRenderedPageBreak markers
# are preserved verbatim where the source document carried them.
$paragraphXml = @(
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''Not the Sun'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''Brand New'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=16813476.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p2.music.126.net/6Ur8Mk5-BweGCZzZLVRc_A==/109951167637374655.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: ''https://redatom.top/lrc/Not the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sun.lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''It''s You'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''Animal Collective'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=16493900.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p2.music.126.net/muuLitiDCJdAEVnbfSBCaQ==/109951165166485770.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: ''''https://redatom.top/lrc/It''s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>You.lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''午後'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''death''s dynamic shroud'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=1392514991.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p1.music.126.net/LOOWgydP-ZdCL5zWC0mn9w==/109951164379853882.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: ''''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''The Other Side'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''The Dismemberment Plan'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=19416116.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p1.music.126.net/5FQJEAhf7vEkTrnYtrdzMA==/109951166350043916.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: ''https://redatom.top/lrc/The Other </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Side.lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''Liberation Frequency'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''Refused'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=18751260.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p2.music.126.net/LDjG47cLPkMG41txhOHW8w==/109951164497608984.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: ''https://redatom.top/lrc/Liberation </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frequency.lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''2:35 (Version 2)'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''Spacemen 3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=3430981.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">            cover: ''http://p2.music.126.net/ZXtenTd_P13JKQZ3YDlv6g==/1777910302120005.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: ''''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            name: ''Call of the Wild'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            artist: ''SB The Moor'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            url: ''http://music.163.com/song/media/outer/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>url?id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=1831494316.mp3'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            cover: ''http://p1.music.126.net/zgKJM7fZaANkRfS-oLHeOQ==/109951165825714963.jpg'',</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">            </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lrc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: ''''</w:t></w:r></w:p>',
    '<w:p><w:r><w:t xml:space="preserve">        },</w:t></w:r></w:p>'
)

$d = $word.ActiveDocument

# The document starts with a single empty paragraph; replace it with the first
# fragment, then keep targeting the next untouched paragraph for each subsequent one.
for ($i = 0; $i -lt $paragraphXml.Length; $i++) {
    $target = $d.Paragraphs($i + 1).Range
    $target.InsertXML($paragraphXml[$i])
}
